$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (预测值) updates
$ws.Range("C2").Value = -570.5
$ws.Range("C3").Value = -574.2
$ws.Range("C4").Value = -694.2
$ws.Range("C7").Value = -751.3
$ws.Range("C8").Value = -714.8
$ws.Range("C9").Value = -651.4
$ws.Range("C13").Value = 134.6
$ws.Range("C14").Value = 183.5
$ws.Range("C15").Value = 76.5
$ws.Range("C18").Value = 100.3
$ws.Range("C19").Value = 256
$ws.Range("C22").Value = -111.7
$ws.Range("C23").Value = -51.8

# Column B (真实值) update
$ws.Range("B3").Value = -841.3
